$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "·规则：理工农医类,社会调查报告和人文社科类等3个选项`n·样例：社会调查报告和人文社科类"
$ws.Range("E6").Value = "·规则：哲学,经济学,法学,教育学,文学等13个选项`n·样例：经济学"
$ws.Range("F6").Value = "·规则：哲学,理论经济学,应用经济学,法学等111个选项`n·样例：测绘科学与技术"
$ws.Range("L6").Value = "·规则：数学科学学院,物理学院等73个选项`n·样例：汇丰商学院"
$ws.Range("R6").Value = "·规则：男,女`n·样例：男"
$ws.Range("S6").Value = "·规则：数学科学学院,物理学院等73个选项`n·样例：环境与能源学院"
